$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2.166745185852051

# Row 3
$ws.Range("B3").Value = 0.8879855713787864
$ws.Range("C3").Value = 0.4092276933397176
$ws.Range("D3").Value = 28.95592632604821
$ws.Range("E3").Value = 1734.705056190491

# Row 4
$ws.Range("B4").Value = 10.00378516870295
$ws.Range("C4").Value = 1.195313983374676
$ws.Range("D4").Value = 983209881637749.8
$ws.Range("E4").Value = 0.6918625831604004

# Row 5
$ws.Range("A5").Value = "XGB_MANY"
$ws.Range("B5").Value = 10.92465838353673
$ws.Range("C5").Value = 1.217024147963954
$ws.Range("D5").Value = 185650793619020.4
$ws.Range("E5").Value = 0.4392876625061035

# Row 6
$ws.Range("A6").Value = "LSTM_MANY"
$ws.Range("B6").Value = 1.057351772670686
$ws.Range("C6").Value = 0.4376544579442828
$ws.Range("D6").Value = 6.363815463289402
$ws.Range("E6").Value = 1863.676489830017

# Row 7
$ws.Range("A7").Value = "FOREST_MANY"
$ws.Range("B7").Value = 11.36471627530992
$ws.Range("C7").Value = 1.24376069214876
$ws.Range("D7").Value = 117645374315565.2
$ws.Range("E7").Value = 0.8285977840423584
